# "generar litigio (mediación o juicio)" — update the Ambiente/URL pair used
# for the mediación/juicio claim and bump the NroSiniestro test values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: drop the "i-" prefix from the pre-production host/URL pair.
# (Write the URL/hyperlink cell first so new shared-string entries land in
# the same order the original authoring session produced.)
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"

# Row 3 NroSiniestro: new claim number, kept as text (quote-prefixed, like
# the original cell) with the trailing spaces from the source data.
$ws.Range("F3").Formula = "'1120170200917  "

# Reset the view: scroll back to column A and move the selection to D12.
$ws.Range("D12").Select()
